$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra blank line (row 2) that sits between the header row (row 1)
# and the placeholder data row (originally row 3). This shifts row 3 up to row 2.
$ws.Rows.Item(2).Delete()

# Update the active selection to match the target state.
$ws.Range("D13").Select()
